$d = $word.ActiveDocument

function ReplaceText($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: could not find [$old]"
    }
}

# 1) " + 29,874 + 29,874" -> " + 50 + 50"
ReplaceText "29,874 + 29,874" "50 + 50"

# 2) "39" + "737" + ",88 " (after "p = ") -> "40140,4" + " "
ReplaceText "39737,88 N" "40140,4 N"

# 3) second "39737,88" occurrence (after "  => ") -> "40140,4"; also drops en-US from lang, handled separately below
ReplaceText "  => 39737,88 x" "  => 40140,4 x"

# 4) " P = 2037" + ",84 Pa " -> " P = " + "2058,48" + " Pa "
ReplaceText "19500  P = 2037,84 Pa" "19500  P = 2058,48 Pa"

# 5) "2037" + ",84 " (standalone, start of new paragraph) -> "2058,48" + " "
ReplaceText "2037,84 - 2007,2" "2058,48 - 2007,2"

# 6) " h2 = 17,97 m" -> " h2 = " + "19,98" + " m"
ReplaceText "(h2 -15)  h2 = 17,97 m" "(h2 -15)  h2 = 19,98 m"

# 7) "h2 - h1 => 17,97 - 15 => 2,97 m " -> "h2 - h1 => " + "19,98" + " - 15 => " + "4,98" + " m "
ReplaceText "h2 – h1 => 17,97 – 15 => 2,97 m" "h2 – h1 => 19,98 – 15 => 4,98 m"

# 8) "container ship, when we place 200 containers it will sink 2.97 m." ->
#    "...00 containers it will sink 4,98 m."
ReplaceText "when we place 200 containers it will sink 2.97 m." "when we place 200 containers it will sink 4,98 m."
